# Auto-generated edit script applying the Asura_Profits diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1000
$ws.Range("I52").Value = 1000
$ws.Range("J52").Value = 1000
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 3000
$ws.Range("M52").Value = -2840
$ws.Range("N52").Value = -3320
$ws.Range("H74").Value = 3746.3845
$ws.Range("I74").Value = 3380.6
$ws.Range("J74").Value = 3975
$ws.Range("K74").Value = 3380.6
$ws.Range("L74").Value = 3975
$ws.Range("M74").Value = -2444.6
$ws.Range("N74").Value = -5847
$ws.Range("H77").Value = 3746.3845
$ws.Range("I77").Value = 3380.6
$ws.Range("J77").Value = 3975
$ws.Range("K77").Value = 16903
$ws.Range("L77").Value = 19875
$ws.Range("M77").Value = -12223
$ws.Range("N77").Value = -29235
$ws.Range("H88").Value = 2286.4285
$ws.Range("I88").Value = 1153
$ws.Range("J88").Value = 2475.3333
$ws.Range("K88").Value = 1153
$ws.Range("L88").Value = 2475.3333
$ws.Range("M88").Value = -747
$ws.Range("N88").Value = -3287.3333
$ws.Range("H91").Value = 2286.4285
$ws.Range("I91").Value = 1153
$ws.Range("J91").Value = 2475.3333
$ws.Range("K91").Value = 1153
$ws.Range("L91").Value = 2475.3333
$ws.Range("M91").Value = 251
$ws.Range("N91").Value = -5283.3333
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 112946
$ws.Range("J134").Value = 112946
$ws.Range("L134").Value = 112946
$ws.Range("N134").Value = -123086
$ws.Range("H135").Value = 1266.5714
$ws.Range("I135").Value = 1145.091
$ws.Range("J135").Value = 1712
$ws.Range("K135").Value = 10305.819
$ws.Range("L135").Value = 15408
$ws.Range("M135").Value = -7770.819
$ws.Range("N135").Value = -20478
$ws.Range("H138").Value = 4361.94
$ws.Range("I138").Value = 2444.228
$ws.Range("J138").Value = 6904.0234
$ws.Range("K138").Value = 7332.684
$ws.Range("L138").Value = 20712.0702
$ws.Range("M138").Value = -2192.684
$ws.Range("N138").Value = -30992.0702

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 35214
$ws.Range("J123").Value = 35214
$ws.Range("L123").Value = 35214
$ws.Range("N123").Value = -45014
$ws.Range("H132").Value = 1554.9073
$ws.Range("I132").Value = 1270.1459
$ws.Range("J132").Value = 3833
$ws.Range("K132").Value = 3810.4377
$ws.Range("L132").Value = 11499
$ws.Range("M132").Value = -1280.4377
$ws.Range("N132").Value = -16559
$ws.Range("H133").Value = 30462.111
$ws.Range("J133").Value = 30462.111
$ws.Range("L133").Value = 30462.111
$ws.Range("N133").Value = -35522.111
$ws.Range("H141").Value = 50914
$ws.Range("J141").Value = 50914
$ws.Range("L141").Value = 50914
$ws.Range("N141").Value = -61274

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 18000
$ws.Range("J2").Value = 18000
$ws.Range("L2").Value = 18000
$ws.Range("N2").Value = -18226

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2315.8333
$ws.Range("J16").Value = 2947.5
$ws.Range("L16").Value = 2947.5
$ws.Range("N16").Value = -3521.5
$ws.Range("H31").Value = 2383.302
$ws.Range("I31").Value = 1393.1666
$ws.Range("J31").Value = 6163.8184
$ws.Range("K31").Value = 1393.1666
$ws.Range("L31").Value = 6163.8184
$ws.Range("M31").Value = -1098.1666
$ws.Range("N31").Value = -6753.8184
$ws.Range("H34").Value = 2383.302
$ws.Range("I34").Value = 1393.1666
$ws.Range("J34").Value = 6163.8184
$ws.Range("K34").Value = 1393.1666
$ws.Range("L34").Value = 6163.8184
$ws.Range("M34").Value = -1191.1666
$ws.Range("N34").Value = -6567.8184
$ws.Range("H105").Value = 5314.6523
$ws.Range("I105").Value = 7351.8
$ws.Range("J105").Value = 1495
$ws.Range("K105").Value = 7351.8
$ws.Range("L105").Value = 1495
$ws.Range("M105").Value = -5604.8
$ws.Range("N105").Value = -4989
$ws.Range("H107").Value = 744.44446
$ws.Range("I107").Value = 587.5
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 587.5
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1332.5
$ws.Range("N107").Value = -5840
$ws.Range("H113").Value = 2315.8333
$ws.Range("J113").Value = 2947.5
$ws.Range("L113").Value = 2947.5
$ws.Range("N113").Value = -7287.5
$ws.Range("H132").Value = 1450.5714
$ws.Range("I132").Value = 1256.081
$ws.Range("J132").Value = 2889.8
$ws.Range("K132").Value = 3768.242999999999
$ws.Range("L132").Value = 8669.400000000001
$ws.Range("M132").Value = -1238.242999999999
$ws.Range("N132").Value = -13729.4
$ws.Range("H135").Value = 135430
$ws.Range("J135").Value = 135430
$ws.Range("L135").Value = 135430
$ws.Range("N135").Value = -145570
$ws.Range("H137").Value = 49158.57
$ws.Range("J137").Value = 74703.336
$ws.Range("L137").Value = 74703.336
$ws.Range("N137").Value = -84903.336
$ws.Range("H140").Value = 76126
$ws.Range("J140").Value = 76126
$ws.Range("L140").Value = 76126
$ws.Range("N140").Value = -86486

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 11799.8
$ws.Range("I70").Value = 11799.8
$ws.Range("K70").Value = 35399.39999999999
$ws.Range("M70").Value = -35084.39999999999
$ws.Range("H73").Value = 11799.8
$ws.Range("I73").Value = 11799.8
$ws.Range("K73").Value = 35399.39999999999
$ws.Range("M73").Value = -34307.39999999999
$ws.Range("H87").Value = 10000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 10000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H112").Value = 3692.8386
$ws.Range("I112").Value = 1190.6
$ws.Range("J112").Value = 4174.0386
$ws.Range("K112").Value = 3571.8
$ws.Range("L112").Value = 12522.1158
$ws.Range("M112").Value = -2463.8
$ws.Range("N112").Value = -14738.1158

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4335
$ws.Range("I80").Value = 4002.5
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 4002.5
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -3004.5
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 4335
$ws.Range("I83").Value = 4002.5
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 20012.5
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -15020.5
$ws.Range("N83").Value = -34984
$ws.Range("H132").Value = 2263.1177
$ws.Range("I132").Value = 1676.7858
$ws.Range("K132").Value = 5030.357400000001
$ws.Range("M132").Value = -2500.357400000001
$ws.Range("H133").Value = 54103.332
$ws.Range("J133").Value = 54103.332
$ws.Range("L133").Value = 54103.332
$ws.Range("N133").Value = -64223.332

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1376.375
$ws.Range("I16").Value = 1376.375
$ws.Range("K16").Value = 1376.375
$ws.Range("M16").Value = -1206.375
$ws.Range("H55").Value = 352
$ws.Range("I55").Value = 352
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 352
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -179
$ws.Range("N55").ClearContents()
$ws.Range("H82").Value = 1354.1538
$ws.Range("J82").Value = 1614.2858
$ws.Range("L82").Value = 1614.2858
$ws.Range("N82").Value = -2336.2858
$ws.Range("H85").Value = 1354.1538
$ws.Range("J85").Value = 1614.2858
$ws.Range("L85").Value = 1614.2858
$ws.Range("N85").Value = -4110.2858
$ws.Range("H99").Value = 191253
$ws.Range("I99").Value = 224103.6
$ws.Range("J99").Value = 27000
$ws.Range("K99").Value = 224103.6
$ws.Range("L99").Value = 27000
$ws.Range("M99").Value = -221108.6
$ws.Range("N99").Value = -32990
$ws.Range("H134").Value = 60216
$ws.Range("J134").Value = 60216
$ws.Range("L134").Value = 60216
$ws.Range("N134").Value = -70356
$ws.Range("H138").Value = 56666.668
$ws.Range("J138").Value = 56666.668
$ws.Range("L138").Value = 56666.668
$ws.Range("N138").Value = -66946.66800000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1793.7894
$ws.Range("I136").Value = 1810.1111
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 5430.3333
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -2880.3333
$ws.Range("N136").Value = -9600
$ws.Range("H138").Value = 120428
$ws.Range("J138").Value = 120428
$ws.Range("L138").Value = 120428
$ws.Range("N138").Value = -130708
